$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds version-number-like text (e.g. "311.98", "27.953.20").
# Plain `.Value = "..."` assignment lets Excel auto-coerce numeric-looking
# strings into actual numbers, which would flip the stored cell type from
# string to number. Force text storage by setting the cell format to Text
# ("@") before writing the literal, then restore the original "Normal" style
# (General format, no borders) so cell styling matches the source file.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.953.20"
$ws.Range("E2").Value = "  -0.47%  "
Set-TextValue "D3" "1.857.37"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "311.98"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.07%  "
Set-TextValue "D7" "0.5130"
$ws.Range("E7").Value = "  +1.71%  "
Set-TextValue "D8" "0.3813"
$ws.Range("E8").Value = "  -0.55%  "
Set-TextValue "D9" "0.08225"
$ws.Range("E9").Value = "  -5.13%  "
$ws.Range("E10").Value = "  -0.79%  "
Set-TextValue "D11" "41.46"
$ws.Range("E11").Value = "  -0.09%  "
Set-TextValue "D12" "6.170"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "20.46"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.848.29"
$ws.Range("E14").Value = "  -1.27%  "
Set-TextValue "D15" "7.248"
$ws.Range("E16").Value = "  -0.05%  "
Set-TextValue "D17" "0.00001094"
$ws.Range("E17").Value = "  -0.76%  "
Set-TextValue "D18" "90.28"
$ws.Range("E18").Value = "  -0.77%  "
Set-TextValue "D19" "0.06641"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  -0.08%  "
Set-TextValue "D22" "6.006"
$ws.Range("E22").Value = "  -1.53%  "
Set-TextValue "D23" "27.993.21"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  -3.61%  "
Set-TextValue "D25" "2.240"
$ws.Range("E25").Value = "  -1.14%  "
Set-TextValue "D26" "2.071.96"
$ws.Range("E26").Value = "  -0.80%  "
Set-TextValue "D27" "2.501"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -1.42%  "
Set-TextValue "D30" "124.42"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("E33").Value = "  +4.46%  "
Set-TextValue "D34" "3.595"
$ws.Range("E34").Value = "  -0.08%  "
Set-TextValue "D35" "9.380"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("E36").Value = "  -1.74%  "
Set-TextValue "D37" "0.06492"
$ws.Range("E37").Value = "  -1.29%  "
Set-TextValue "D38" "0.2178"
$ws.Range("E38").Value = "  +0.18%  "
Set-TextValue "D39" "0.6558"
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("E40").Value = "  -1.04%  "
Set-TextValue "D41" "4.974"
$ws.Range("E41").Value = "  +1.64%  "
Set-TextValue "D42" "1.207"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("E43").Value = "  -3.87%  "
Set-TextValue "D44" "0.6106"
$ws.Range("E44").Value = "  +2.09%  "
Set-TextValue "D45" "13.02"
Set-TextValue "D46" "3.678"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  -0.66%  "
Set-TextValue "D48" "2.006"
$ws.Range("E48").Value = "  +0.99%  "
Set-TextValue "D49" "1.210"
$ws.Range("E49").Value = "  -1.55%  "
Set-TextValue "D50" "120.67"
$ws.Range("E50").Value = "  -0.43%  "
Set-TextValue "D51" "77.93"
$ws.Range("E51").Value = "  -2.62%  "
